# feat(indonesian): change language to indonesian
#
# Translate the employee-template header row (row 1, columns A:N) from
# English to Indonesian. Each column keeps the same meaning/position —
# only the shared-string text changes. Also restores the previously
# recorded selection (E7) that accompanied the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nama Lengkap"       # was: Fullname
$ws.Range("B1").Value = "NIP"                # was: Employee ID
$ws.Range("C1").Value = "NIK"                # unchanged
$ws.Range("D1").Value = "NPWP"               # unchanged
$ws.Range("E1").Value = "Agama"              # was: Religion
$ws.Range("F1").Value = "Tempat Lahir"       # was: Place of Birth
$ws.Range("G1").Value = "Nomor Handphone"    # was: Phone Number
$ws.Range("H1").Value = "Jenis Kelamin"      # was: Gender
$ws.Range("I1").Value = "Status Pernikahan"  # was: Marital Status
$ws.Range("J1").Value = "Golongan Darah"     # was: Blood Type
$ws.Range("K1").Value = "Email"              # unchanged
$ws.Range("L1").Value = "Password Akun"      # was: Account Password
$ws.Range("M1").Value = "Alamat"             # was: Address
$ws.Range("N1").Value = "Alamat KTP"         # was: Citizen Address

# Restore the recorded cursor position/selection for the sheet (scrolled so
# column D is the left-most visible column, with E7 as the active cell).
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("E7").Select()
